# Refresh the cryptocurrency price/volume snapshot (Coin, Link, Price,
# Volume(1h) columns) to match the latest scrape, including two rank swaps
# (Uniswap/BitcoinCash at rows 19-20, and BabyDogeCoin/Quant at rows 43-44).
#
# Note: several "Price" values look numeric (e.g. "1.001", "0.06349") but
# must stay literal text - exactly like the original cells - otherwise
# Excel would silently convert them to numbers and drop significant
# trailing/formatting digits. Prefixing those with a leading apostrophe
# forces Excel to store them as text (the same thing happens if you type
# an apostrophe before a number in the Excel UI), while leaving values
# that are unambiguous as text (URLs, names, multi-dot prices, percentages
# with surrounding spaces, subscript-digit prices) alone.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.901.89"
$ws.Range("E2").Value = "  -0.20%  "
$ws.Range("D3").Value = "1.634.69"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "  -0.47%  "
$ws.Range("D5").Value = "'213.92"
$ws.Range("E5").Value = "  -0.79%  "
$ws.Range("D6").Value = "'0.5056"
$ws.Range("E6").Value = "  -0.53%  "
$ws.Range("D7").Value = "'1.001"
$ws.Range("E7").Value = "  -0.35%  "
$ws.Range("D8").Value = "'0.2566"
$ws.Range("E8").Value = "  +0.18%  "
$ws.Range("D9").Value = "'0.06349"
$ws.Range("E9").Value = "  -0.52%  "
$ws.Range("D10").Value = "'19.65"
$ws.Range("E10").Value = "  +1.00%  "
$ws.Range("D11").Value = "'0.07745"
$ws.Range("E11").Value = "  -0.46%  "
$ws.Range("D12").Value = "'4.279"
$ws.Range("E12").Value = "  -0.45%  "
$ws.Range("D13").Value = "1.633.46"
$ws.Range("E13").Value = "  -0.84%  "
$ws.Range("D14").Value = "'0.5430"
$ws.Range("E14").Value = "  -0.67%  "
$ws.Range("D15").Value = "0.0₅7732"
$ws.Range("E15").Value = "  -1.48%  "
$ws.Range("D16").Value = "'64.00"
$ws.Range("E16").Value = "  -0.70%  "
$ws.Range("D17").Value = "25.914.00"
$ws.Range("E17").Value = "  -0.42%  "
$ws.Range("E18").Value = "  -0.33%  "
$ws.Range("B19").Value = "Uniswap"
$ws.Range("C19").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D19").Value = "'4.429"
$ws.Range("E19").Value = "  -0.43%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'194.92"
$ws.Range("E20").Value = "  -1.58%  "
$ws.Range("D21").Value = "'9.906"
$ws.Range("E21").Value = "  -0.60%  "
$ws.Range("D22").Value = "'6.107"
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("E23").Value = "  -0.40%  "
$ws.Range("D24").Value = "'1.894"
$ws.Range("E24").Value = "  +0.82%  "
$ws.Range("D25").Value = "'142.90"
$ws.Range("E25").Value = "  +1.37%  "
$ws.Range("E26").Value = "  +8.20%  "
$ws.Range("D27").Value = "'6.810"
$ws.Range("E27").Value = "  -1.05%  "
$ws.Range("D28").Value = "'15.60"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("D29").Value = "'1.235"
$ws.Range("E29").Value = "  -0.49%  "
$ws.Range("D30").Value = "'0.04874"
$ws.Range("E30").Value = "  -3.05%  "
$ws.Range("E31").Value = "  -0.80%  "
$ws.Range("D32").Value = "'3.193"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("D33").Value = "'1.544"
$ws.Range("E33").Value = "  +0.11%  "
$ws.Range("E34").Value = "  +0.12%  "
$ws.Range("D35").Value = "'0.9089"
$ws.Range("E35").Value = "  +1.15%  "
$ws.Range("D36").Value = "'2.570"
$ws.Range("E36").Value = "  -1.10%  "
$ws.Range("D37").Value = "'0.5495"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "1.122.65"
$ws.Range("E38").Value = "  -0.95%  "
$ws.Range("D39").Value = "'0.01559"
$ws.Range("E39").Value = "  -0.11%  "
$ws.Range("E40").Value = "  -0.42%  "
$ws.Range("D41").Value = "'5.579"
$ws.Range("E41").Value = "  -0.71%  "
$ws.Range("D42").Value = "'0.8037"
$ws.Range("E42").Value = "  -1.85%  "
$ws.Range("B43").Value = "BabyDogeCoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D43").Value = "0.0₈124"
$ws.Range("E43").Value = "  -8.74%  "
$ws.Range("B44").Value = "Quant"
$ws.Range("C44").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D44").Value = "'98.47"
$ws.Range("E44").Value = "  -1.84%  "
$ws.Range("D45").Value = "1.768.68"
$ws.Range("E45").Value = "  -0.64%  "
$ws.Range("D46").Value = "'0.4479"
$ws.Range("E46").Value = "  -1.29%  "
$ws.Range("D47").Value = "'1.004"
$ws.Range("E47").Value = "  +0.20%  "
$ws.Range("D48").Value = "'54.94"
$ws.Range("E48").Value = "  -0.08%  "
$ws.Range("D49").Value = "'0.05172"
$ws.Range("E49").Value = "  +1.95%  "
$ws.Range("D50").Value = "'7.487"
$ws.Range("E50").Value = "  +1.18%  "
$ws.Range("D51").Value = "'1.001"
$ws.Range("E51").Value = "  -0.54%  "
